# Restore the "adductName" column (as column C) on the "Corrected" sheet of
# the accucor workbook. This shifts the existing sample-intensity columns
# (previously C:F) one column to the right (now D:G).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Corrected")

# Insert a new blank column at C, pushing C:F -> D:G.
$ws.Columns("C:C").Insert()

# The freshly inserted column inherits bold+centered formatting from the
# column to its left (B). Reset it to the base "Normal" style first so that
# the bold-only header style we build next doesn't also carry the centered
# alignment.
$ws.Range("C1").Style = "Normal"

# Header cell: "adductName" (bold, left/general aligned - distinct from the
# other bold+centered headers).
$ws.Range("C1").Value = "adductName"
$ws.Range("C1").Font.Bold = $true

# Data rows: the adduct-name/label-ratio values for each compound+label row.
$ws.Range("C2").Value = 1.0402279999999999
$ws.Range("C3").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("C6").Value = 1.078981
$ws.Range("C7").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("C10").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("C12").Value = 0
